# Auto-generated Excel COM-interop script
# Applies market-price data refresh updates to the Jenova_Profits workbook
# (columns H-N: currentAveragePrice*, LevePrice*, LeveProfit* cached values)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 886.4375
$ws.Range("I6").Value = 1186.1818
$ws.Range("J6").Value = 227
$ws.Range("K6").Value = 3558.5454
$ws.Range("L6").Value = 681
$ws.Range("M6").Value = -3446.5454
$ws.Range("N6").Value = -905
$ws.Range("H17").Value = 2166.6538
$ws.Range("J17").Value = 2166.6538
$ws.Range("L17").Value = 6499.9614
$ws.Range("N17").Value = -6835.9614
$ws.Range("H62").Value = 15631595
$ws.Range("J62").Value = 7998.3335
$ws.Range("L62").Value = 7998.3335
$ws.Range("N62").Value = -9246.333500000001
$ws.Range("H65").Value = 15631595
$ws.Range("J65").Value = 7998.3335
$ws.Range("L65").Value = 39991.6675
$ws.Range("N65").Value = -46231.6675
$ws.Range("H100").Value = 10481.412
$ws.Range("I100").Value = 3366
$ws.Range("J100").Value = 14362.546
$ws.Range("K100").Value = 3366
$ws.Range("L100").Value = 14362.546
$ws.Range("M100").Value = -2825
$ws.Range("N100").Value = -15444.546
$ws.Range("H132").Value = 1866.6586
$ws.Range("I132").Value = 1705.8206
$ws.Range("K132").Value = 5117.4618
$ws.Range("M132").Value = -2587.4618
$ws.Range("H138").Value = 4515.2104
$ws.Range("I138").Value = 3401.2
$ws.Range("J138").Value = 4622.327
$ws.Range("K138").Value = 10203.6
$ws.Range("L138").Value = 13866.981
$ws.Range("M138").Value = -5063.599999999999
$ws.Range("N138").Value = -24146.981

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4053.9487
$ws.Range("I32").Value = 3777.6575
$ws.Range("K32").Value = 3777.6575
$ws.Range("M32").Value = -3490.6575
$ws.Range("H44").Value = 43274.5
$ws.Range("J44").Value = 43274.5
$ws.Range("L44").Value = 43274.5
$ws.Range("N44").Value = -44250.5
$ws.Range("H45").Value = 3212.1538
$ws.Range("I45").Value = 2157.5
$ws.Range("K45").Value = 2157.5
$ws.Range("M45").Value = -1780.5
$ws.Range("H74").Value = 1421.6428
$ws.Range("I74").Value = 1360.1052
$ws.Range("K74").Value = 1360.1052
$ws.Range("M74").Value = -486.1052
$ws.Range("H77").Value = 1421.6428
$ws.Range("I77").Value = 1360.1052
$ws.Range("K77").Value = 6800.526
$ws.Range("M77").Value = -2432.526
$ws.Range("H80").Value = 86037.5
$ws.Range("J80").Value = 89716.664
$ws.Range("L80").Value = 89716.664
$ws.Range("N80").Value = -91712.664
$ws.Range("H83").Value = 86037.5
$ws.Range("J83").Value = 89716.664
$ws.Range("L83").Value = 269149.992
$ws.Range("N83").Value = -279133.992
$ws.Range("H102").Value = 1925.7142
$ws.Range("I102").Value = 1976.9412
$ws.Range("J102").Value = 1708
$ws.Range("K102").Value = 1976.9412
$ws.Range("L102").Value = 1708
$ws.Range("M102").Value = -354.9412
$ws.Range("N102").Value = -4952
$ws.Range("H132").Value = 2124.6482
$ws.Range("I132").Value = 2173.196
$ws.Range("K132").Value = 6519.588
$ws.Range("M132").Value = -3989.588

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 770.6667
$ws.Range("I64").Value = 600.6667
$ws.Range("J64").Value = 813.1667
$ws.Range("K64").Value = 600.6667
$ws.Range("L64").Value = 813.1667
$ws.Range("M64").Value = -375.6667
$ws.Range("N64").Value = -1263.1667
$ws.Range("H67").Value = 770.6667
$ws.Range("I67").Value = 600.6667
$ws.Range("J67").Value = 813.1667
$ws.Range("K67").Value = 600.6667
$ws.Range("L67").Value = 813.1667
$ws.Range("M67").Value = 179.3333
$ws.Range("N67").Value = -2373.1667
$ws.Range("H86").Value = 896674.9
$ws.Range("I86").Value = 1309324.6
$ws.Range("K86").Value = 1309324.6
$ws.Range("M86").Value = -1308201.6
$ws.Range("H89").Value = 896674.9
$ws.Range("I89").Value = 1309324.6
$ws.Range("K89").Value = 6546623
$ws.Range("M89").Value = -6541007
$ws.Range("H132").Value = 50000
$ws.Range("J132").Value = 50000
$ws.Range("L132").Value = 50000
$ws.Range("N132").Value = -60120
$ws.Range("H134").Value = 56724.367
$ws.Range("I134").Value = 4213.4375
$ws.Range("J134").Value = 336782.66
$ws.Range("K134").Value = 12640.3125
$ws.Range("L134").Value = 1010347.98
$ws.Range("M134").Value = -10105.3125
$ws.Range("N134").Value = -1015417.98

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 52094.668
$ws.Range("I31").Value = 1000
$ws.Range("J31").Value = 54649.4
$ws.Range("K31").Value = 1000
$ws.Range("L31").Value = 54649.4
$ws.Range("M31").Value = -705
$ws.Range("N31").Value = -55239.4
$ws.Range("H34").Value = 52094.668
$ws.Range("I34").Value = 1000
$ws.Range("J34").Value = 54649.4
$ws.Range("K34").Value = 1000
$ws.Range("L34").Value = 54649.4
$ws.Range("M34").Value = -798
$ws.Range("N34").Value = -55053.4
$ws.Range("H94").Value = 1483.6
$ws.Range("I94").Value = 308.5
$ws.Range("K94").Value = 308.5
$ws.Range("M94").Value = 142.5
$ws.Range("H132").Value = 1440.28
$ws.Range("I132").Value = 1362.9546
$ws.Range("J132").Value = 2007.3334
$ws.Range("K132").Value = 4088.8638
$ws.Range("L132").Value = 6022.0002
$ws.Range("M132").Value = -1558.8638
$ws.Range("N132").Value = -11082.0002
$ws.Range("H133").Value = 34999.5
$ws.Range("I133").Value = 34999.5
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 34999.5
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -32469.5
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 353.85715
$ws.Range("J2").Value = 469.4
$ws.Range("L2").Value = 2816.4
$ws.Range("N2").Value = -3042.4
$ws.Range("H15").Value = 205.5
$ws.Range("I15").Value = 200
$ws.Range("J15").Value = 206.6
$ws.Range("K15").Value = 600
$ws.Range("L15").Value = 619.8
$ws.Range("M15").Value = -460
$ws.Range("N15").Value = -899.8
$ws.Range("H69").Value = 1000
$ws.Range("I69").Value = 1000
$ws.Range("K69").Value = 3000
$ws.Range("M69").Value = -2189
$ws.Range("H72").Value = 1000
$ws.Range("I72").Value = 1000
$ws.Range("K72").Value = 9000
$ws.Range("M72").Value = -4944
$ws.Range("H101").Value = 12543
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 12543
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 37629
$ws.Range("M101").ClearContents()
$ws.Range("N101").Value = -42497

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1195.3
$ws.Range("I97").Value = 1195.3
$ws.Range("K97").Value = 1195.3
$ws.Range("M97").Value = -699.3
$ws.Range("H108").Value = 15000
$ws.Range("J108").Value = 15000
$ws.Range("L108").Value = 15000
$ws.Range("N108").Value = -22680
$ws.Range("H132").Value = 59496.684
$ws.Range("I132").Value = 5965.125
$ws.Range("K132").Value = 17895.375
$ws.Range("M132").Value = -15365.375
$ws.Range("H135").Value = 250100000
$ws.Range("J135").Value = 250100000
$ws.Range("L135").Value = 250100000
$ws.Range("N135").Value = -250110140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6689.5557
$ws.Range("I7").Value = 7246.5454
$ws.Range("J7").Value = 5814.2856
$ws.Range("K7").Value = 7246.5454
$ws.Range("L7").Value = 5814.2856
$ws.Range("M7").Value = -7134.5454
$ws.Range("N7").Value = -6038.2856
$ws.Range("H16").Value = 632.5714
$ws.Range("I16").Value = 671.3333
$ws.Range("J16").Value = 400
$ws.Range("K16").Value = 671.3333
$ws.Range("L16").Value = 400
$ws.Range("M16").Value = -501.3333
$ws.Range("N16").Value = -740
$ws.Range("H46").Value = 4584.8
$ws.Range("I46").Value = 3946.4
$ws.Range("K46").Value = 3946.4
$ws.Range("M46").Value = -3758.4
$ws.Range("H126").Value = 6689.5557
$ws.Range("I126").Value = 7246.5454
$ws.Range("J126").Value = 5814.2856
$ws.Range("K126").Value = 21739.6362
$ws.Range("L126").Value = 17442.8568
$ws.Range("M126").Value = -19269.6362
$ws.Range("N126").Value = -22382.8568
$ws.Range("H131").Value = 60000
$ws.Range("J131").Value = 60000
$ws.Range("L131").Value = 60000
$ws.Range("N131").Value = -70080

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 701.3333
$ws.Range("I126").Value = 1304
$ws.Range("K126").Value = 3912
$ws.Range("M126").Value = -1442
